# "Generate Report for Handoff"
#
# The localization engine produced a fresh handoff for
# a484e22f-441f-430c-8cf1-cb945e5b75c0.md (row 3 on every sheet): its
# status flips from "In Translation" to "Ready for handoff", a new
# machine-translation ("mt") priority is recorded, and new handoff
# timestamps / xliff-generation timestamps are written. Widening the
# status columns is Excel's own column auto-fit reacting to the new,
# longer "Ready for handoff" text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet: row 3 is the a484e22f-441f-430c-8cf1-cb945e5b75c0.md file ---
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-09-06 18:19:21"

# --- zh-cn sheet: row 3 is the a484e22f-441f-430c-8cf1-cb945e5b75c0.md file ---
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("H3").Value = "2016-09-06 18:19:15"

# --- de-de sheet: row 3 is the a484e22f-441f-430c-8cf1-cb945e5b75c0.md file ---
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("H3").Value = "2016-09-06 18:19:21"

# --- Column widths widen to fit the new, longer "Ready for handoff" status text ---
$wsOverview.Columns.Item(5).ColumnWidth = 16.3333333333333
$wsOverview.Columns.Item(6).ColumnWidth = 16.3333333333333
$wsZhCn.Columns.Item(3).ColumnWidth = 16.3333333333333
$wsDeDe.Columns.Item(3).ColumnWidth = 16.3333333333333
